# DOMA-4173: add a new "order" column to the payments export template.
#
# The sheet currently has 8 data columns (A..H): date, account, address,
# unitName, type, transaction, status, amount (for both the i18n header
# row and the two templated rows below it). We need to insert a new
# "order" column between "transaction" (F) and "status" (old G), so the
# final layout (A..I) becomes:
#   date, account, address, unitName, type, transaction, order, status, amount
#
# We avoid Range/Column .Insert() here because on this engine it (a)
# drops the inherited default column style ("style" attribute) on the
# shifted <col> entries and (b) off-by-one's the trailing catch-all
# <col> range's max index. Instead we manually shift the two rightmost
# columns' cell content/formatting one column to the right (bounded to
# the used rows, so it stays fast) and then populate the freed-up column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift old column H (status-row formatting/values) into new column I,
# then old column G (status/amount placeholders' neighbour) into H.
# Doing H before G avoids clobbering G before it is copied onward.
$ws.Range("H1:H10").Copy($ws.Range("I1:I10"))
$ws.Range("G1:G10").Copy($ws.Range("H1:H10"))

# Populate the freed column G with the new "order" placeholders.
$ws.Range("G1").Value = "{d.i18n.order}"
$ws.Range("G2").Value = "{d.objs[I].order}"
$ws.Range("G3").Value = "{d.objs[I+1].order}"

# Widen the new column G (~25.94 chars wide in the target file). This
# engine quantises ColumnWidth to 1/6-character steps on save, so an
# input of 25.17 is the closest achievable round-trip to 25.9375.
$ws.Columns.Item(7).ColumnWidth = 25.17
